# Apply updated Leve profit calculations to the Sheets workbook (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 19462.982
$ws.Range("J17").Value = 19462.982
$ws.Range("L17").Value = 58388.946
$ws.Range("N17").Value = -58724.946

$ws.Range("H28").Value = 1134.238
$ws.Range("I28").Value = 1509.2222
$ws.Range("J28").Value = 853
$ws.Range("K28").Value = 1509.2222
$ws.Range("L28").Value = 853
$ws.Range("M28").Value = -1024.2222
$ws.Range("N28").Value = -1823

$ws.Range("H40").Value = 2083.3809
$ws.Range("I40").Value = 2038.5385
$ws.Range("K40").Value = 2038.5385
$ws.Range("M40").Value = -1863.5385

$ws.Range("H41").Value = 428
$ws.Range("I41").Value = 200
$ws.Range("J41").Value = 478.66666
$ws.Range("K41").Value = 200
$ws.Range("L41").Value = 478.66666
$ws.Range("M41").Value = 240
$ws.Range("N41").Value = -1358.66666

$ws.Range("H44").Value = 35000
$ws.Range("J44").Value = 35000
$ws.Range("L44").Value = 35000
$ws.Range("N44").Value = -35924

$ws.Range("H53").Value = 193.11765
$ws.Range("I53").Value = 115.166664
$ws.Range("J53").Value = 235.63637
$ws.Range("K53").Value = 115.166664
$ws.Range("L53").Value = 235.63637
$ws.Range("M53").Value = 521.833336
$ws.Range("N53").Value = -1509.63637

$ws.Range("H62").Value = 1485.125
$ws.Range("I62").Value = 1374.421
$ws.Range("J62").Value = 1905.8
$ws.Range("K62").Value = 1374.421
$ws.Range("L62").Value = 1905.8
$ws.Range("M62").Value = -750.421
$ws.Range("N62").Value = -3153.8

$ws.Range("H65").Value = 1485.125
$ws.Range("I65").Value = 1374.421
$ws.Range("J65").Value = 1905.8
$ws.Range("K65").Value = 6872.105
$ws.Range("L65").Value = 9529
$ws.Range("M65").Value = -3752.105
$ws.Range("N65").Value = -15769

$ws.Range("H74").Value = 4454.5
$ws.Range("I74").Value = 6312.625
$ws.Range("J74").Value = 3392.7144
$ws.Range("K74").Value = 6312.625
$ws.Range("L74").Value = 3392.7144
$ws.Range("M74").Value = -5376.625
$ws.Range("N74").Value = -5264.7144

$ws.Range("H76").Value = 2936.9524
$ws.Range("I76").Value = 2687.4119
$ws.Range("J76").Value = 3997.5
$ws.Range("K76").Value = 2687.4119
$ws.Range("L76").Value = 3997.5
$ws.Range("M76").Value = -2372.4119
$ws.Range("N76").Value = -4627.5

$ws.Range("H77").Value = 4454.5
$ws.Range("I77").Value = 6312.625
$ws.Range("J77").Value = 3392.7144
$ws.Range("K77").Value = 31563.125
$ws.Range("L77").Value = 16963.572
$ws.Range("M77").Value = -26883.125
$ws.Range("N77").Value = -26323.572

$ws.Range("H79").Value = 2936.9524
$ws.Range("I79").Value = 2687.4119
$ws.Range("J79").Value = 3997.5
$ws.Range("K79").Value = 2687.4119
$ws.Range("L79").Value = 3997.5
$ws.Range("M79").Value = -1595.4119
$ws.Range("N79").Value = -6181.5

$ws.Range("H86").Value = 2488.889
$ws.Range("I86").Value = 2066.6667
$ws.Range("J86").Value = 2700
$ws.Range("K86").Value = 2066.6667
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -943.6667000000002
$ws.Range("N86").Value = -4946

$ws.Range("H89").Value = 2488.889
$ws.Range("I89").Value = 2066.6667
$ws.Range("J89").Value = 2700
$ws.Range("K89").Value = 10333.3335
$ws.Range("L89").Value = 13500
$ws.Range("M89").Value = -4717.333500000001
$ws.Range("N89").Value = -24732

$ws.Range("H92").Value = 251.55
$ws.Range("I92").Value = 259.5263
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 259.5263
$ws.Range("L92").Value = 100
$ws.Range("M92").Value = 988.4737
$ws.Range("N92").Value = -2596

$ws.Range("H107").Value = 758.1875
$ws.Range("I107").Value = 556.5
$ws.Range("J107").Value = 1094.3334
$ws.Range("K107").Value = 556.5
$ws.Range("L107").Value = 1094.3334
$ws.Range("M107").Value = 1363.5
$ws.Range("N107").Value = -4934.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10240.154
$ws.Range("I32").Value = 6010.6143
$ws.Range("J32").Value = 35315.285
$ws.Range("K32").Value = 6010.6143
$ws.Range("L32").Value = 35315.285
$ws.Range("M32").Value = -5723.6143
$ws.Range("N32").Value = -35889.285

$ws.Range("H63").Value = 1994.0625
$ws.Range("I63").Value = 1782.2727
$ws.Range("K63").Value = 1782.2727
$ws.Range("M63").Value = -1096.2727

$ws.Range("H66").Value = 1994.0625
$ws.Range("I66").Value = 1782.2727
$ws.Range("K66").Value = 8911.363499999999
$ws.Range("M66").Value = -5479.363499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 361.27274
$ws.Range("I22").Value = 361.27274
$ws.Range("K22").Value = 361.27274
$ws.Range("M22").Value = -188.27274

$ws.Range("H105").Value = 2527365.8
$ws.Range("I105").Value = 3789532
$ws.Range("J105").Value = 3033.3333
$ws.Range("K105").Value = 3789532
$ws.Range("L105").Value = 3033.3333
$ws.Range("M105").Value = -3787785
$ws.Range("N105").Value = -6527.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 250
$ws.Range("I59").Value = 250
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 750
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -210
$ws.Range("N59").ClearContents()

$ws.Range("H70").Value = 2911.25
$ws.Range("I70").Value = 1286.6666
$ws.Range("K70").Value = 3859.9998
$ws.Range("M70").Value = -3544.9998

$ws.Range("H73").Value = 2911.25
$ws.Range("I73").Value = 1286.6666
$ws.Range("K73").Value = 3859.9998
$ws.Range("M73").Value = -2767.9998

$ws.Range("H113").Value = 7936959
$ws.Range("I113").Value = 505.8
$ws.Range("J113").Value = 10989441
$ws.Range("K113").Value = 1517.4
$ws.Range("L113").Value = 32968323
$ws.Range("M113").Value = 652.5999999999999
$ws.Range("N113").Value = -32972663

$ws.Range("H131").Value = 63252.062
$ws.Range("I131").Value = 240.66667
$ws.Range("J131").Value = 87908.69500000001
$ws.Range("K131").Value = 722.00001
$ws.Range("L131").Value = 263726.085
$ws.Range("M131").Value = 4317.99999
$ws.Range("N131").Value = -273806.085

$ws.Range("H136").Value = 597.6
$ws.Range("I136").Value = 540
$ws.Range("J136").Value = 600
$ws.Range("K136").Value = 1620
$ws.Range("L136").Value = 1800
$ws.Range("M136").Value = 3480
$ws.Range("N136").Value = -12000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4345
$ws.Range("I70").Value = 4388.25
$ws.Range("J70").Value = 4310.4
$ws.Range("K70").Value = 4388.25
$ws.Range("L70").Value = 4310.4
$ws.Range("M70").Value = -4118.25
$ws.Range("N70").Value = -4850.4

$ws.Range("H73").Value = 4345
$ws.Range("I73").Value = 4388.25
$ws.Range("J73").Value = 4310.4
$ws.Range("K73").Value = 4388.25
$ws.Range("L73").Value = 4310.4
$ws.Range("M73").Value = -3452.25
$ws.Range("N73").Value = -6182.4

$ws.Range("H80").Value = 64925.25
$ws.Range("I80").Value = 2266.5
$ws.Range("J80").Value = 102520.5
$ws.Range("K80").Value = 2266.5
$ws.Range("L80").Value = 102520.5
$ws.Range("M80").Value = -1268.5
$ws.Range("N80").Value = -104516.5

$ws.Range("H83").Value = 64925.25
$ws.Range("I83").Value = 2266.5
$ws.Range("J83").Value = 102520.5
$ws.Range("K83").Value = 11332.5
$ws.Range("L83").Value = 512602.5
$ws.Range("M83").Value = -6340.5
$ws.Range("N83").Value = -522586.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H132").Value = 2817.348
$ws.Range("I132").Value = 1899.3077
$ws.Range("J132").Value = 4010.8
$ws.Range("K132").Value = 5697.9231
$ws.Range("L132").Value = 12032.4
$ws.Range("M132").Value = -3167.9231
$ws.Range("N132").Value = -17092.4
